# Update TC_CategoryPage (sheet2) column L "Actual" result strings to new concise
# action-log format (e.g. "Element verified succesfully" -> "VerifyElement: null").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC_CategoryPage")

$ws.Range("L2").Value = 'VerifyTitle: City Market Norwalk - Online Grocery Supermarket with Home Delivery'
$ws.Range("L3").Value = 'VerifyElement: null'
$ws.Range("L4").Value = 'Click: null'
$ws.Range("L5").Value = 'VerifyElement: null'
$ws.Range("L6").Value = 'SetText: Randomemailid'
$ws.Range("L7").Value = 'SetText: 123456'
$ws.Range("L8").Value = 'Click: null'
$ws.Range("L9").Value = 'VerifyText: Akash sangal'
$ws.Range("L10").Value = 'VerifyTitle: City Market Norwalk - Online Grocery Supermarket with Home Delivery'
$ws.Range("L11").Value = 'MoveToCategory: Quick & Easy Food Solutions'
$ws.Range("L12").Value = 'VerifyElement: null'
$ws.Range("L13").Value = 'VerifyElement: null'
$ws.Range("L14").Value = 'VerifyElement: null'
$ws.Range("L15").Value = 'VerifyElement: null'
$ws.Range("L16").Value = 'VerifyElement: null'
$ws.Range("L17").Value = 'VerifyElement: null'
$ws.Range("L18").Value = 'VerifyElement: null'
$ws.Range("L19").Value = @"
VerifyText: 717 West Ave Norwalk, Connecticut,
CT 06850 , USA
+1 203-956-0241
"@
$ws.Range("L20").Value = 'VerifyText: All Days 7:30 AM – 9:00 PM'
$ws.Range("L21").Value = 'VerifyText: Copyright © 2018 City Market Norwalk. All rights reserved. Terms Of Use & Privacy Policy'
$ws.Range("L22").Value = 'Click: null'
$ws.Range("L23").Value = 'VerifyTitle: City Market Norwalk - Terms of Use'
$ws.Range("L24").Value = 'CloseBrowser: Child'
$ws.Range("L25").Value = 'Click: null'
$ws.Range("L26").Value = 'VerifyTitle: City Super Market Norwalk - Privacy Policy'
$ws.Range("L27").Value = 'CloseBrowser: Child'
$ws.Range("L28").Value = 'Click: null'
$ws.Range("L29").Value = 'VerifyElement: null'
$ws.Range("L30").Value = 'Click: null'
$ws.Range("L31").Value = 'VerifyNoElement: null'
$ws.Range("L32").Value = 'Click: null'
$ws.Range("L33").Value = 'VerifyTitle: Offers'
$ws.Range("L34").Value = 'MoveToCategory: Quick & Easy Food Solutions'
$ws.Range("L35").Value = 'VerifyTitle: Order online for Quick & Easy Food Solutions, delivery or store pickup|City Market Norwalk'
$ws.Range("L36").Value = 'Click: null'
$ws.Range("L37").Value = 'VerifyElement: null'
$ws.Range("L38").Value = 'HeaderMenuSearch: My Account'
$ws.Range("L39").Value = 'HeaderMenuSearch: Account Information'
$ws.Range("L40").Value = 'HeaderMenuSearch: My Orders'
$ws.Range("L41").Value = 'HeaderMenuSearch: My Saved Cards'
$ws.Range("L42").Value = 'HeaderMenuSearch: My Wish List'
$ws.Range("L43").Value = 'HeaderMenuSearch: Delivery Coverage'
$ws.Range("L44").Value = 'HeaderMenuSearch: Offers'
$ws.Range("L45").Value = 'HeaderMenuSearch: Sign Out'
$ws.Range("L46").Value = 'Click: null'
$ws.Range("L47").Value = 'VerifyNoElement: null'
$ws.Range("L48").Value = 'Click: null'
$ws.Range("L49").Value = 'VerifyElement: null'
$ws.Range("L50").Value = 'VerifyText: You have no items in your shopping cart.'
$ws.Range("L51").Value = 'Click: null'
$ws.Range("L52").Value = 'VerifyNoElement: null'
$ws.Range("L53").Value = 'SetText: milk'
$ws.Range("L54").Value = 'Click: null'
$ws.Range("L55").Value = 'SearchProduct: Skim Plus 100% Fat Free Milk'
$ws.Range("L56").Value = 'MoveToCategory: Quick & Easy Food Solutions'
$ws.Range("L57").Value = 'SetText: asdasd@'
$ws.Range("L58").Value = 'Click: null'
$ws.Range("L59").Value = 'VerifyText: Please enter a valid email address (Ex: johndoe@domain.com).'
$ws.Range("L60").Value = 'SetText: asdasd@assdcsadsaasd.com'
$ws.Range("L61").Value = 'Click: null'
$ws.Range("L62").Value = 'VerifyText: Thank you for your subscription.'
$ws.Range("L63").Value = 'VerifyCategoryAndProduct: null'
$ws.Range("L64").Value = 'VerifySiteMapFromCategoryPage: null'
$ws.Range("L65").Value = 'VerifyFooterLinks: null'
$ws.Range("L66").Value = 'Click: null'
$ws.Range("L67").Value = 'Click: null'
$ws.Range("L68").Value = 'Wait: 6000'
$ws.Range("L69").Value = 'VerifyElement: null'
